$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 3 (Testmail #20) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Kun je deze taak op je nemen?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #20: Kun je deze taak op je nemen?"
$logs.Range("D3").Value = "Overig"
$logs.Range("E3").Value = "Beste [Naam],`nBedankt voor je bericht. Kun je wat meer details geven over welke taak je precies bedoelt? Dan kan ik je zo goed mogelijk helpen.`nMet vriendelijke groet,`n[Jouw naam]  `nE-mailassistent bij [Bedrijfsnaam]"
$logs.Range("F3").Value = "2025-06-27 00:03:01"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Nee"
$logs.Range("I3").Value = "Ja"

# Extend the conditional-formatting blocks down to include row 3
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D3"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G3"))
$logs.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H3"))
$logs.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I3"))

# --- Dashboard sheet: append row 3 (Overig / 1) ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A3").Value = "Overig"
$dash.Range("B3").Value = 1

# --- Chart: extend category/value series references to include row 3 ---
$chart = $dash.ChartObjects(1).Chart
$chart.SeriesCollection(1).XValues = "='Dashboard'!`$A`$2:`$A`$3"
$chart.SeriesCollection(1).Values = "='Dashboard'!`$B`$2:`$B`$3"
